$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update sample number value from E7760 to E7420 (shared string used in G2:G19)
$ws.Range("G2:G19").Value = "E7420"

# 2. Replace boolean H2:H19 cells with a FALSE() formula (keeps same displayed value/result).
#    Set cell-by-cell so Excel does not collapse the column into a single shared formula.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# 3. Update the active selection to G2:G19 (matches the new selection in the diff)
$ws.Range("G2:G19").Select()
